# Generate Report for Handoff
#
# Updates the localization-status report after a fresh handoff run:
#   - Priority for the three still-pending files (28bcd9bb, 4b6f174d,
#     65454548, c7ca2161) flips from "low" to "ht" on both the zh-cn and
#     de-de sheets.
#   - The "Latest Handoff Datetime" for those same rows is refreshed to the
#     new handoff timestamp (different per target locale).
#   - The Overview sheet's "Latest HO Xliff Generate Date" tracks the de-de
#     handoff timestamp, so it is refreshed alongside de-de's column.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$rows = @(4, 5, 6, 7)

foreach ($r in $rows) {
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-08-19 14:36:54"

    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-08-19 14:36:59"

    $overview.Range("G$r").Value = "2016-08-19 14:36:59"
}
